$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for "Espinaca" (Primera, 08-Nov-2021) was added to
# the dataset. In the source data it sits right before the existing
# 09-Jun-2021 record, so insert a fresh row at 225 (pushing everything from
# 225..279 down to 226..280) and fill it with the new observation.
$ws.Rows.Item(225).Insert()

$ws.Range("A225").Value = 9
$ws.Range("B225").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C225").Value = "Metropolitana"
$ws.Range("D225").Value = 44508
$ws.Range("E225").Value = 13
$ws.Range("F225").Value = 100112012
$ws.Range("G225").Value = "Espinaca"
$ws.Range("H225").Value = "Sin especificar"
$ws.Range("I225").Value = "Primera"
$ws.Range("J225").Value = 124
$ws.Range("K225").Value = 5000
$ws.Range("L225").Value = 6000
$ws.Range("M225").Value = 5500
$ws.Range("N225").Value = "`$/cuna 10 kilos"
$ws.Range("O225").Value = "Provincia de Chacabuco"
$ws.Range("P225").Value = 550
$ws.Range("Q225").Value = 10
$ws.Range("R225").Value = "Hortaliza"
